$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: change date text format
$ws.Range("A2").Value = "26/08/2025"

# D2: change from text "001984" to numeric 1984
$ws.Range("D2").Value = 1984

# H2: change from text "6255258016" to numeric 6255258016
$ws.Range("H2").Value = 6255258016

# I2: remove comma after "Añez"
$ws.Range("I2").Value = "En La Esquina Del Hospital Municipal De Cotoca Por La Entrada De Emergencia Rodolfo Añez Nro.: S/N Zona/Barrio.: COTOCA"
